$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is numeric-looking need NumberFormat forced to
# Text first so Excel keeps them as literal strings (matching the source
# data, which stores these as text, e.g. "1.000", "0.9989").

$ws.Range("D2").Value = "27.964.57"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "1.772.64"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.44"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9989"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4499"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07435"
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.95"
$ws.Range("E10").Value = "  +1.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.107"
$ws.Range("E11").Value = "  +2.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9994"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.99"
$ws.Range("E13").Value = "  +2.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.046"
$ws.Range("E14").Value = "  +2.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.251"
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("D16").Value = "1.768.88"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.71"
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001065"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06434"
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9985"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.16"
$ws.Range("E21").Value = "  +2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.784"
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("D23").Value = "27.982.11"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.53"
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("E27").Value = "  +1.65%  "
$ws.Range("D28").Value = "1.971.84"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.163"
$ws.Range("E29").Value = "  +6.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.68"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.110"
$ws.Range("E31").Value = "  +6.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.697"
$ws.Range("E32").Value = "  +6.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09215"
$ws.Range("E33").Value = "  +1.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.690"
$ws.Range("E34").Value = "  +1.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.87"
$ws.Range("E35").Value = "  +2.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06215"
$ws.Range("E36").Value = "  +3.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02292"
$ws.Range("E37").Value = "  +1.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2111"
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.978"
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6318"
$ws.Range("E40").Value = "  +1.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.184"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.397"
$ws.Range("E42").Value = "  +1.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.904"
$ws.Range("E43").Value = "  +2.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.35"
$ws.Range("E44").Value = "  +1.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.751"
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5882"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.69"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.961"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06890"
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.90"
$ws.Range("E51").Value = "  +2.53%  "
